# Refresh the crypto price/volume table (GitHub Actions scheduled update).
# Note: several "Price" column values look numeric (e.g. "417.60", "1.00") but
# must stay plain text, matching the workbook's original inlineStr cells
# (trailing zeros / digit grouping would otherwise be lost to numeric
# coercion). A leading apostrophe forces Excel to store them as text, exactly
# like typing '417.60 into a cell - the apostrophe itself is not stored.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.892.80'
$ws.Range("E2").Value = '  +6.48%  '
$ws.Range("D3").Value = '3.558.68'
$ws.Range("E3").Value = '  +2.57%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").Value = '''417.60'
$ws.Range("E5").Value = '  +0.40%  '
$ws.Range("D6").Value = '''129.65'
$ws.Range("E6").Value = '  -0.77%  '
$ws.Range("D7").Value = '''0.650'
$ws.Range("E7").Value = '  +3.91%  '
$ws.Range("D8").Value = '3.549.53'
$ws.Range("E8").Value = '  +2.43%  '
$ws.Range("E9").Value = '  -0.15%  '
$ws.Range("D10").Value = '''0.782'
$ws.Range("E10").Value = '  +7.02%  '
$ws.Range("D11").Value = '''0.181'
$ws.Range("E11").Value = '  +27.14%  '
$ws.Range("D12").Value = '''0.0000351'
$ws.Range("E12").Value = '  +58.87%  '
$ws.Range("D13").Value = '''42.67'
$ws.Range("E13").Value = '  -0.45%  '
$ws.Range("D14").Value = '''9.95'
$ws.Range("E14").Value = '  +2.89%  '
$ws.Range("D15").Value = '4.115.86'
$ws.Range("E15").Value = '  +2.40%  '
$ws.Range("E16").Value = '  -0.25%  '
$ws.Range("D17").Value = '''20.20'
$ws.Range("E17").Value = '  -1.70%  '
$ws.Range("D18").Value = '3.533.52'
$ws.Range("E18").Value = '  +0.92%  '
$ws.Range("D19").Value = '''1.13'
$ws.Range("E19").Value = '  +4.42%  '
$ws.Range("D20").Value = '''12.43'
$ws.Range("E20").Value = '  -3.37%  '
$ws.Range("D21").Value = '66.716.54'
$ws.Range("E21").Value = '  +6.31%  '
$ws.Range("D22").Value = '''456.68'
$ws.Range("E22").Value = '  -3.96%  '
$ws.Range("D23").Value = '''89.82'
$ws.Range("E23").Value = '  -1.37%  '
$ws.Range("E24").Value = '  -3.33%  '
$ws.Range("D25").Value = '''13.03'
$ws.Range("E25").Value = '  -3.40%  '
$ws.Range("D26").Value = '''3.37'
$ws.Range("E26").Value = '  +1.45%  '
$ws.Range("D27").Value = '''9.95'
$ws.Range("E27").Value = '  -5.74%  '
$ws.Range("D28").Value = '''34.70'
$ws.Range("E28").Value = '  +3.89%  '
$ws.Range("D29").Value = '''4.84'
$ws.Range("E29").Value = '  +0.60%  '
$ws.Range("D30").Value = '''2.78'
$ws.Range("E30").Value = '  +5.12%  '
$ws.Range("D31").Value = '''12.41'
$ws.Range("E31").Value = '  +2.83%  '
$ws.Range("D32").Value = '''0.117'
$ws.Range("E32").Value = '  +3.66%  '
$ws.Range("D33").Value = '''7.27'
$ws.Range("E33").Value = '  -4.47%  '
$ws.Range("D34").Value = '''0.159'
$ws.Range("E34").Value = '  -5.22%  '
$ws.Range("D35").Value = '''1.00'
$ws.Range("E35").Value = '  -0.06%  '
$ws.Range("D36").Value = '''39.16'
$ws.Range("E36").Value = '  -4.76%  '
$ws.Range("D37").Value = '''56.73'
$ws.Range("E37").Value = '  -2.65%  '
$ws.Range("B38").Value = 'PEPE'
$ws.Range("C38").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D38").Value = '0.0₃0786'
$ws.Range("E38").Value = '  +37.94%  '
$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").Value = '''0.0494'
$ws.Range("E39").Value = '  +0.32%  '
$ws.Range("D40").Value = '''0.148'
$ws.Range("E40").Value = '  +9.99%  '
$ws.Range("E41").Value = '  +0.01%  '
$ws.Range("B42").Value = 'Monero'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D42").Value = '''148.37'
$ws.Range("E42").Value = '  +2.20%  '
$ws.Range("B43").Value = 'WEMIXToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D43").Value = '''2.72'
$ws.Range("E43").Value = '  +0.68%  '
$ws.Range("B44").Value = 'Stacks'
$ws.Range("C44").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D44").Value = '''2.95'
$ws.Range("E44").Value = '  -3.37%  '
$ws.Range("D45").Value = '''4.35'
$ws.Range("E45").Value = '  -1.45%  '
$ws.Range("D46").Value = '''3.23'
$ws.Range("E46").Value = '  -4.05%  '
$ws.Range("E47").Value = '  -5.40%  '
$ws.Range("E48").Value = '  -5.66%  '
$ws.Range("D49").Value = '''2.29'
$ws.Range("E49").Value = '  -4.41%  '
$ws.Range("B50").Value = 'BitcoinSV'
$ws.Range("C50").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range("D50").Value = '''118.42'
$ws.Range("E50").Value = '  +7.59%  '
$ws.Range("B51").Value = 'ApeXProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D51").Value = '''2.59'
$ws.Range("E51").Value = '  +10.49%  '
